$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.455.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.49%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.21'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E5").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5160'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.81%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3942'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07734'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("E9").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E10").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.10'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.34%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.287'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.582'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.001'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.826.93'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.36'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.78%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001082'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.54%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06622'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.072'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.62%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.462.06'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.257'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.04%  '
$ws.Range("E25").ClearFormats()
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("B26").ClearFormats()
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C26").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.57'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E26").ClearFormats()
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("B27").ClearFormats()
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C27").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.449'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.74%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.64'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.034.07'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.01'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.131'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.661'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.664'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.37%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07161'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2238'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.97%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.974'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.69%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02328'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.161'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.57%  '
$ws.Range("E39").ClearFormats()
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Aptos'
$ws.Range("B40").ClearFormats()
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C40").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.28'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("E40").ClearFormats()
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("B41").ClearFormats()
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C41").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6249'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.75%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.191'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("E42").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.395'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.49'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.32%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5907'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.09%  '
$ws.Range("E46").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.982'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.27%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.183'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06923'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.08%  '
$ws.Range("E51").ClearFormats()
